$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-shuffled the per-row market data (date, volume,
# prices, packaging unit, price/kg and kg-per-unit, quality) across rows
# 2-18. Column layout: A Mercado ID, B Mercado, C Region, D Fecha,
# E Codreg, F Categoria ID, G Categoria, H Variedad, I Calidad,
# J Volumen, K Precio minimo, L Precio maximo, M Precio promedio,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion.

$rows = @{
    2  = @{ D = 44313; I = "Primera"; J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 }
    3  = @{ D = 44313; I = "Primera"; J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    4  = @{ D = 45041; I = "Primera"; J = 80; K = 24000; L = 24000; M = 24000; N = "`$/caja 18 kilos empedrada"; P = 1333; Q = 18 }
    5  = @{ D = 45040; I = "Primera"; J = 80; K = 18000; L = 18000; M = 18000; N = "`$/caja 18 kilos empedrada"; P = 1000; Q = 18 }
    6  = @{ D = 45034; I = "Primera"; J = 50; K = 24000; L = 24000; M = 24000; N = "`$/caja 18 kilos granel";    P = 1333; Q = 18 }
    7  = @{ D = 45014; I = "Primera"; J = 30; K = 8000;  L = 8000;  M = 8000;  N = "`$/caja 18 kilos empedrada"; P = 444;  Q = 18 }
    8  = @{ D = 44280; I = "Primera"; J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    10 = @{ D = 45042; I = "Primera"; J = 60; K = 24000; L = 24000; M = 24000; N = "`$/caja 18 kilos empedrada"; P = 1333; Q = 18 }
    11 = @{ D = 44285; I = "Primera"; J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    12 = @{ D = 44293; I = "Primera"; J = 10; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos empedrada"; P = 1667; Q = 15 }
    13 = @{ D = 45044; I = "Primera"; J = 40; K = 24000; L = 24000; M = 24000; N = "`$/caja 18 kilos empedrada"; P = 1333; Q = 18 }
    14 = @{ D = 44315; I = "Especial"; J = 10; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    15 = @{ D = 44315; I = "Primera"; J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos granel";    P = 1000; Q = 15 }
    17 = @{ D = 45037; I = "Primera"; J = 80; K = 24000; L = 24000; M = 24000; N = "`$/caja 15 kilos empedrada"; P = 1600; Q = 15 }
    18 = @{ D = 45015; I = "Primera"; J = 50; K = 24000; L = 24000; M = 24000; N = "`$/caja 18 kilos empedrada"; P = 1333; Q = 18 }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 4).Value2 = $row.D    # D: Fecha
    $ws.Cells.Item($r, 9).Value = $row.I     # I: Calidad
    $ws.Cells.Item($r, 10).Value2 = $row.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value2 = $row.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $row.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $row.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row.N    # N: Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value2 = $row.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value2 = $row.Q   # Q: Kg o Unidades
}
